$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$rows = @(4, 16, 17, 21, 24, 25)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "x"
}

$ws.Range("A22").Select()

$wb.Save()
